# Auto-generated edit script: update Kraken_Profits market-price snapshot values
# across multiple worksheets, per the scheduled-runner diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H86").Value = 2935.875
$ws.Range("I86").Value = 2958.4
$ws.Range("J86").Value = 2898.3333
$ws.Range("K86").Value = 2958.4
$ws.Range("L86").Value = 2898.3333
$ws.Range("M86").Value = -1835.4
$ws.Range("N86").Value = -5144.3333
$ws.Range("H89").Value = 2935.875
$ws.Range("I89").Value = 2958.4
$ws.Range("J89").Value = 2898.3333
$ws.Range("K89").Value = 14792
$ws.Range("L89").Value = 14491.6665
$ws.Range("M89").Value = -9176
$ws.Range("N89").Value = -25723.6665
$ws.Range("H106").Value = 3995
$ws.Range("I106").Value = 3995
$ws.Range("K106").Value = 3995
$ws.Range("M106").Value = -3364
$ws.Range("H111").Value = 295
$ws.Range("I111").Value = 295
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 885
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 2182
$ws.Range("N111").ClearContents()
$ws.Range("H132").Value = 4305.4814
$ws.Range("I132").Value = 3463.1
$ws.Range("K132").Value = 10389.3
$ws.Range("M132").Value = -7859.299999999999
$ws.Range("H135").Value = 2045
$ws.Range("J135").Value = 3000
$ws.Range("L135").Value = 27000
$ws.Range("N135").Value = -32070
$ws.Range("H137").Value = 4873.875
$ws.Range("I137").Value = 4415.3335
$ws.Range("K137").Value = 13246.0005
$ws.Range("M137").Value = -10696.0005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2865.25
$ws.Range("J74").Value = 4599.8
$ws.Range("L74").Value = 4599.8
$ws.Range("N74").Value = -6347.8
$ws.Range("H77").Value = 2865.25
$ws.Range("J77").Value = 4599.8
$ws.Range("L77").Value = 22999
$ws.Range("N77").Value = -31735

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 890.125
$ws.Range("I80").Value = 777.25
$ws.Range("J80").Value = 1003
$ws.Range("K80").Value = 777.25
$ws.Range("L80").Value = 1003
$ws.Range("M80").Value = 220.75
$ws.Range("N80").Value = -2999
$ws.Range("H83").Value = 890.125
$ws.Range("I83").Value = 777.25
$ws.Range("J83").Value = 1003
$ws.Range("K83").Value = 3886.25
$ws.Range("L83").Value = 5015
$ws.Range("M83").Value = 1105.75
$ws.Range("N83").Value = -14999
$ws.Range("H105").Value = 4252.5
$ws.Range("I105").Value = 3804
$ws.Range("K105").Value = 3804
$ws.Range("M105").Value = -2057

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1534.625
$ws.Range("I16").Value = 1396.1666
$ws.Range("J16").Value = 1950
$ws.Range("K16").Value = 1396.1666
$ws.Range("L16").Value = 1950
$ws.Range("M16").Value = -1109.1666
$ws.Range("N16").Value = -2524
$ws.Range("H31").Value = 5990.8667
$ws.Range("I31").Value = 3871.2856
$ws.Range("J31").Value = 7845.5
$ws.Range("K31").Value = 3871.2856
$ws.Range("L31").Value = 7845.5
$ws.Range("M31").Value = -3576.2856
$ws.Range("N31").Value = -8435.5
$ws.Range("H34").Value = 5990.8667
$ws.Range("I34").Value = 3871.2856
$ws.Range("J34").Value = 7845.5
$ws.Range("K34").Value = 3871.2856
$ws.Range("L34").Value = 7845.5
$ws.Range("M34").Value = -3669.2856
$ws.Range("N34").Value = -8249.5
$ws.Range("H60").Value = 21000
$ws.Range("I60").Value = 8333.333000000001
$ws.Range("J60").Value = 25222.223
$ws.Range("K60").Value = 8333.333000000001
$ws.Range("L60").Value = 25222.223
$ws.Range("M60").Value = -7822.333000000001
$ws.Range("N60").Value = -26244.223
$ws.Range("H86").Value = 3500
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 3500
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H113").Value = 1534.625
$ws.Range("I113").Value = 1396.1666
$ws.Range("J113").Value = 1950
$ws.Range("K113").Value = 1396.1666
$ws.Range("L113").Value = 1950
$ws.Range("M113").Value = 773.8334
$ws.Range("N113").Value = -6290

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2625000
$ws.Range("I11").Value = 2857142.8
$ws.Range("K11").Value = 2857142.8
$ws.Range("M11").Value = -2857003.8
$ws.Range("H80").Value = 2399.8
$ws.Range("J80").Value = 1750
$ws.Range("L80").Value = 1750
$ws.Range("N80").Value = -3746
$ws.Range("H83").Value = 2399.8
$ws.Range("J83").Value = 1750
$ws.Range("L83").Value = 8750
$ws.Range("N83").Value = -18734
$ws.Range("H107").Value = 92.72727
$ws.Range("I107").Value = 92
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 92
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = 1828
$ws.Range("N107").Value = -3940
$ws.Range("H113").Value = 1405.4
$ws.Range("I113").Value = 1405.4
$ws.Range("K113").Value = 1405.4
$ws.Range("M113").Value = 764.5999999999999
$ws.Range("H132").Value = 4185.625
$ws.Range("I132").Value = 4185.625
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12556.875
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10026.875
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 24038
$ws.Range("I50").Value = 24038
$ws.Range("K50").Value = 24038
$ws.Range("M50").Value = -23401
$ws.Range("H61").Value = 5822
$ws.Range("I61").Value = 5233.5
$ws.Range("J61").Value = 6999
$ws.Range("K61").Value = 5233.5
$ws.Range("L61").Value = 6999
$ws.Range("M61").Value = -5031.5
$ws.Range("N61").Value = -7403
$ws.Range("H100").Value = 2334.5
$ws.Range("I100").Value = 2432.7778
$ws.Range("J100").Value = 1450
$ws.Range("K100").Value = 2432.7778
$ws.Range("L100").Value = 1450
$ws.Range("M100").Value = -1891.7778
$ws.Range("N100").Value = -2532
$ws.Range("H113").Value = 5822
$ws.Range("I113").Value = 5233.5
$ws.Range("J113").Value = 6999
$ws.Range("K113").Value = 5233.5
$ws.Range("L113").Value = 6999
$ws.Range("M113").Value = -3063.5
$ws.Range("N113").Value = -11339

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3106.6428
$ws.Range("I122").Value = 2899.4
$ws.Range("J122").Value = 3624.75
$ws.Range("K122").Value = 8698.200000000001
$ws.Range("L122").Value = 10874.25
$ws.Range("M122").Value = -6248.200000000001
$ws.Range("N122").Value = -15774.25
